$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.60753999999999997
$ws.Range("B3").Value = 0.69030999999999998
$ws.Range("B4").Value = 0.59492
$ws.Range("B5").Value = 0.59145999999999999
$ws.Range("B6").Value = 0.48859000000000002
$ws.Range("B7").Value = 0.40660000000000002
$ws.Range("B8").Value = 0.38080999999999998
$ws.Range("B9").Value = 0.36770000000000003
$ws.Range("B10").Value = 0.40527000000000002
$ws.Range("B11").Value = 0.18786
$ws.Range("B12").Value = 0.34936
$ws.Range("B13").Value = 0.21362999999999999

$ws.Range("G11").Select()
